$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-12) updated values
$ws.Range("C2").Value = 6.428615036886185
$ws.Range("C3").Value = 1.179444761626655
$ws.Range("C4").Value = 72.19119401741773
$ws.Range("C5").Value = 0.001212478087836644
$ws.Range("C6").Value = 6800.642091352325
$ws.Range("C7").Value = 225.860721036559
$ws.Range("C8").Value = 12.55734109721379
$ws.Range("C9").Value = 0.04257072621840052
$ws.Range("C10").Value = 93.12300139479339
$ws.Range("C11").Value = 0.6570096939249197
$ws.Range("C12").Value = 22.64693908020854

# Column B (rows 13-22) updated values
$ws.Range("B13").Value = 2.526329814456403
$ws.Range("B14").Value = 0.01740844821324572
$ws.Range("B15").Value = 5440.07830106176
$ws.Range("B16").Value = 0.0003606278914958239
$ws.Range("B17").Value = 66.62428904755507
$ws.Range("B18").Value = 45.05167570587946
$ws.Range("B19").Value = 2.40918464056449
$ws.Range("B20").Value = 0.01204602731013438
$ws.Range("B21").Value = 0.6570096939249197
$ws.Range("B22").Value = 22.16662792675197
